$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the formatting of the existing
# header cells (e.g. G1 "sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Corresponding data value for row 2
$ws.Range("H2").Value = 0
